# Appointment_Outcomes.xlsx - mark the pharmacist-pending prescriptions
# for P1001/D002 (row 3) and P1002/D001 (row 4) as DISPENSED now that the
# pharmacist has fulfilled them. Column G holds "Prescription Status".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G3").Value = "DISPENSED"
$ws.Range("G4").Value = "DISPENSED"
